{"js": "// The budget section about average office rents in Gy\u0151r reads\n// \"...d\u00edja Gy\u0151rben kb 70 \u2013 300 ezer forint k\u00f6z\u00f6tt...\" \u2014 \"kb\" is short for\n// \"k\u00f6r\u00fclbel\u00fcl\" (approx.) and is missing its abbreviation-dot. Add the\n// missing period so it reads \"...Gy\u0151rben kb. 70 \u2013 300 ezer forint...\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst needle = \"Gy\u0151rben kb 70\";\nconst replacement = \"Gy\u0151rben kb. 70\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.indexOf(needle) !== -1) {\n    const fixedText = paragraph.text.split(needle).join(replacement);\n    paragraph.getRange().insertText(fixedText, \"Replace\");\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "# The budget section about average office rents in Gy\u0151r reads\n# \"...d\u00edja Gy\u0151rben kb 70 - 300 ezer forint k\u00f6z\u00f6tt...\" - \"kb\" is short for\n# \"k\u00f6r\u00fclbel\u00fcl\" (approx.) and is missing its abbreviation-dot. Add the\n# missing period so it reads \"...Gy\u0151rben kb. 70 - 300 ezer forint...\".\n$d = $word.ActiveDocument\n\nforeach ($paragraph in $d.Paragraphs) {\n    $rng = $paragraph.Range\n    if ($rng.Text.Contains(\"Nincs sz\u00fcks\u00e9g\u00fcnk t\u00fal nagy irod\u00e1ra\")) {\n        $find = $rng.Find\n        $find.Text = \"Gy\u0151rben kb 70\"\n        $find.Replacement.Text = \"Gy\u0151rben kb. 70\"\n        $find.Execute([ref]\"Gy\u0151rben kb 70\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"Gy\u0151rben kb. 70\", [ref]1) | Out-Null\n        break\n    }\n}\n"}
